$d = $word.ActiveDocument

$pairs = @(
    @("180÷7=", "463÷5="),
    @("349÷9=", "829÷8="),
    @("714÷3=", "579÷6="),
    @("593÷9=", "702÷4="),
    @("700÷3=", "822÷7="),
    @("332÷6=", "503÷5="),
    @("634÷6=", "854÷9="),
    @("645÷8=", "231÷3="),
    @("771÷3=", "375÷9="),
    @("936÷9=", "712÷3="),
    @("512÷7=", "657÷5="),
    @("239÷2=", "840÷4="),
    @("517÷2=", "658÷8="),
    @("610÷4=", "323÷2="),
    @("635÷9=", "915÷7="),
    @("108÷4=", "455÷5="),
    @("822÷3=", "232÷7="),
    @("987÷9=", "470÷7="),
    @("612÷9=", "463÷9="),
    @("279÷9=", "115÷7="),
    @("521÷6=", "305÷9="),
    @("456÷9=", "724÷7="),
    @("180÷5=", "351÷7="),
    @("324÷5=", "732÷6="),
    @("602÷6=", "328÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
